$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top, pushing all existing data down by one row.
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "Country"
$ws.Range("B1").Value = "Population"

# Match the saved selection/view state from the edit (cursor on B2, no special
# top-left scroll position anymore).
$ws.Range("B2").Select()
